# Generate Report for Handback
# - Update Status text from "Ready for handoff" to "Handed back: in sync with en-US"
# - Refresh "Latest Handback DateTime" for both locales
# - Clear the (now resolved) "Error Detail" column
# - Widen the Status / Latest Handback DateTime columns, narrow the Error Detail column

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Columns E (zh-cn) and F (de-de) grow wider to fit the new status text
$wsOverview.Range("E1").ColumnWidth = 29.15
$wsOverview.Range("F1").ColumnWidth = 29.15

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsZhCn.Range("K2").Value = "2016-08-03 02:31:19"
$wsZhCn.Range("K3").Value = "2016-08-03 02:31:19"

$wsZhCn.Range("P2").Value = ""
$wsZhCn.Range("P3").Value = ""

# Status column widens, Error Detail column shrinks now that it is empty
$wsZhCn.Range("C1").ColumnWidth = 29.15
$wsZhCn.Range("P1").ColumnWidth = 12.86

# --- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

$wsDeDe.Range("K2").Value = "2016-08-03 02:31:35"
$wsDeDe.Range("K3").Value = "2016-08-03 02:31:35"

$wsDeDe.Range("P2").Value = ""
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Range("C1").ColumnWidth = 29.15
$wsDeDe.Range("P1").ColumnWidth = 12.86
